# Город Сокровищ / Марилл.xlsx — add two new rows (29, 30) of translated
# script-line data after the existing table, and convert row 28 into a
# "continuation" row (border style) like the rest of its group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 28 becomes a "continuation" row: gets an explicit (empty) A28
#    cell and its border/font style switches from the "group head"
#    style (4/5) to the "continuation" style (6/7), matching rows such
#    as 5/7/10/14/18/21/23/25.
# ---------------------------------------------------------------------
$ws.Range("A25:B25").Copy() | Out-Null
$ws.Range("A28:B28").PasteSpecial(-4122) | Out-Null
$ws.Range("C25:E25").Copy() | Out-Null
$ws.Range("C28:E28").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Rows.Item(28).RowHeight = 15.6

# ---------------------------------------------------------------------
# 2. New row 29 — a new "group head" row (script file name in column A).
# ---------------------------------------------------------------------
$ws.Range("A29").Value = "SCRIPT/T01P02A/us2301.ssb"
$ws.Range("B29").Value = 19
$ws.Range("C29").Value = ' Oh, [hero] and\n[partner]!'
$ws.Range("D29").Value = " О, [hero] и [partner]!"
$ws.Range("E29").Value = " Ï, [hero] é [partner]!"

$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A29:B29").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:E2").Copy() | Out-Null
$ws.Range("C29:E29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Rows.Item(29).RowHeight = 43.2

# ---------------------------------------------------------------------
# 3. New row 30 — continuation of the us2301 block (no script name).
# ---------------------------------------------------------------------
$ws.Range("B30").Value = 22
$ws.Range("C30").Value = ' Mr. [CS:N]Drowzee[CR] left to go travel a\nwhile ago.'
$ws.Range("D30").Value = ' Не так давно, Мистер [CS:N]Дроузи[CR]\nушёл в путешествие.'
$ws.Range("E30").Value = ' Îå óàë äàâîï, Íéòóåñ [CS:N]Äñïôèé[CR]\nôšæì â ðôóåšåòóâéå.'

$ws.Range("B13").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null
$ws.Range("C13:E13").Copy() | Out-Null
$ws.Range("C30:E30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Rows.Item(30).RowHeight = 31.8

# ---------------------------------------------------------------------
# 4. Scroll/selection bookkeeping to match the saved view state.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("C30").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 28

Write-Host "edit applied"
